# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.151.68'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.900.37'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = "'306.91"
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').Value = "'0.5231"
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('D8').Value = "'0.3805"
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('D9').Value = "'0.07293"
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').Value = "'21.35"
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').Value = "'0.08206"
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').Value = '1.882.87'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').Value = "'95.54"
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = "'5.355"
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = "'14.69"
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = '27.190.72'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = "'5.124"
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').Value = '2.126.44'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = "'10.78"
$ws.Range('E23').Value = '  +1.73%  '
$ws.Range('D24').Value = "'6.461"
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = "'2.325"
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'149.14"
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('D27').Value = "'18.29"
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = "'1.742"
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('D29').Value = "'115.44"
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('D31').Value = "'4.891"
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').Value = "'0.09218"
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').Value = "'0.05045"
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = "'0.7934"
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('D35').Value = "'1.225"
$ws.Range('E35').Value = '  -1.17%  '
$ws.Range('D36').Value = "'2.966"
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').Value = "'3.360"
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').Value = "'2.639"
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('D39').Value = "'0.5739"
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').Value = "'0.01991"
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('E41').Value = '  +0.65%  '
$ws.Range('D42').Value = "'9.093"
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').Value = "'116.51"
$ws.Range('E44').Value = '  -1.58%  '
$ws.Range('D45').Value = "'0.1518"
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('D46').Value = "'0.4910"
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('D47').Value = "'1.002"
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('D48').Value = "'10.18"
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').Value = "'1.635"
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').Value = "'38.54"
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('E51').Value = '  +0.84%  '

Write-Output "Applied 94 cell updates"
